$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.866.14'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("D3").Value = '2.350.21'
$ws.Range("E3").Value = '  +0.59%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = '''0.672'
$ws.Range("E5").Value = '  +3.53%  '

$ws.Range("D6").Value = '''237.84'
$ws.Range("E6").Value = '  +2.63%  '

$ws.Range("D7").Value = '''72.58'
$ws.Range("E7").Value = '  +10.42%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").Value = '  +18.94%  '

$ws.Range("E10").Value = '  +3.95%  '

$ws.Range("D11").Value = '''28.92'
$ws.Range("E11").Value = '  +8.52%  '

$ws.Range("E12").Value = '  +2.59%  '

$ws.Range("D13").Value = '2.699.82'
$ws.Range("E13").Value = '  +0.44%  '

$ws.Range("D14").Value = '''16.75'
$ws.Range("E14").Value = '  +9.38%  '

$ws.Range("E15").Value = '  +7.04%  '

$ws.Range("D16").Value = '''0.898'
$ws.Range("E16").Value = '  +6.59%  '

$ws.Range("D17").Value = '2.358.42'
$ws.Range("E17").Value = '  +0.91%  '

$ws.Range("D18").Value = '43.846.35'
$ws.Range("E18").Value = '  +1.43%  '

$ws.Range("E19").Value = '  +4.14%  '

$ws.Range("D20").Value = '''77.93'
$ws.Range("E20").Value = '  +5.42%  '

$ws.Range("D21").Value = '''6.40'
$ws.Range("E21").Value = '  +3.71%  '

$ws.Range("D22").Value = '''253.80'
$ws.Range("E22").Value = '  +2.24%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("D24").Value = '''3.76'
$ws.Range("E24").Value = '  -3.09%  '

$ws.Range("D25").Value = '''2.50'
$ws.Range("E25").Value = '  +3.44%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''10.50'
$ws.Range("E26").Value = '  +6.26%  '

$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '''2.28'
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").Value = '''22.38'
$ws.Range("E28").Value = '  +1.17%  '

$ws.Range("D29").Value = '''172.68'
$ws.Range("E29").Value = '  -1.30%  '

$ws.Range("E30").Value = '  +6.42%  '

$ws.Range("E31").Value = '  +1.89%  '

$ws.Range("E32").Value = '  +5.22%  '

$ws.Range("D33").Value = '''5.16'
$ws.Range("E33").Value = '  +3.01%  '

$ws.Range("E34").Value = '  +4.40%  '

$ws.Range("D35").Value = '''5.24'
$ws.Range("E35").Value = '  +6.18%  '

$ws.Range("D36").Value = '''3.96'
$ws.Range("E36").Value = '  +10.40%  '

$ws.Range("E37").Value = '  -2.58%  '

$ws.Range("D38").Value = '''6.40'
$ws.Range("E38").Value = '  -0.65%  '

$ws.Range("E39").Value = '  +6.55%  '

$ws.Range("D40").Value = '''19.48'
$ws.Range("E40").Value = '  +8.68%  '

$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("E42").Value = '  -1.69%  '

$ws.Range("E43").Value = '  +4.22%  '

$ws.Range("D44").Value = '''0.0981'
$ws.Range("E44").Value = '  +4.25%  '

$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Value = '''4.45'
$ws.Range("E46").Value = '  +0.92%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''98.06'
$ws.Range("E47").Value = '  -0.57%  '

$ws.Range("D48").Value = '''0.180'
$ws.Range("E48").Value = '  +11.19%  '

$ws.Range("E49").Value = '  +3.86%  '

$ws.Range("D50").Value = '1.436.03'
$ws.Range("E50").Value = '  +0.09%  '

$ws.Range("E51").Value = '  +1.42%  '
